# Populate the previously-blank row 5 of the Places Mapping sheet with a new
# mapping rule: "ca_places.access" (mirrors the pattern used by rows 2-4,
# which each map a "Mapping" rule type / sequence number / CA field target).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Mapping"
$ws.Range("B5").Value = 7
$ws.Range("C5").Value = "ca_places.access"
